# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" (and before "总计"),
#    populated with the fund holdings for 2022-Q1 (same column layout as "2021-Q4").
# 2. Prepend a new summary row for "2022-Q1" to the "总计" sheet (above the
#    existing "2021-Q4" row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the "2022-Q1" sheet, positioned between "2021-Q4" and "总计"
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add($null, $q4)
$q1.Name = "2022-Q1"

# Header row + index column: copy formatting from the "2021-Q4" sheet so the
# new sheet matches its look (bold, centered, bordered header/index cells).
$q4.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("A2:A7").Copy()
$q1.Range("A2:A7").PasteSpecial(-4122)

$q1.Cells.Item(1, 2).Value = "基金代码"
$q1.Cells.Item(1, 3).Value = "基金名称"
$q1.Cells.Item(1, 4).Value = "基金规模"
$q1.Cells.Item(1, 5).Value = "股票总仓位"
$q1.Cells.Item(1, 6).Value = "仓位占比"
$q1.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q1.Cells.Item(1, 8).Value = "仓位排名"

$rows = @(
    @(0, "010490", "鹏华高质量增长混合A",           "13.31", "93.61", "3.94", "0.5244", 4),
    @(1, "009023", "鹏华稳健回报混合",               "3.52",  "93.91", "5.26", "0.1852", 7),
    @(2, "011351", "金鹰年年邮益一年持有期混合A",   "9.03",  "37.02", "0.91", "0.0822", 10),
    @(3, "004265", "金鹰民丰回报定期开放混合",       "6.57",  "28.61", "0.70", "0.0460", 7),
    @(4, "010491", "鹏华高质量增长混合C",           "0.28",  "93.61", "3.94", "0.0110", 4),
    @(5, "011352", "金鹰年年邮益一年持有期混合C",   "0.59",  "37.02", "0.91", "0.0054", 10)
)

$r = 2
foreach ($row in $rows) {
    $q1.Cells.Item($r, 1).Value = $row[0]
    $q1.Cells.Item($r, 2).Value = "'" + $row[1]
    $q1.Cells.Item($r, 3).Value = $row[2]
    $q1.Cells.Item($r, 4).Value = "'" + $row[3]
    $q1.Cells.Item($r, 5).Value = "'" + $row[4]
    $q1.Cells.Item($r, 6).Value = "'" + $row[5]
    $q1.Cells.Item($r, 7).Value = "'" + $row[6]
    $q1.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2) Update the "总计" sheet: insert a new row above the existing data for the
#    new "2022-Q1" summary line (shift "2021-Q4" down from row 2 to row 3).
# ---------------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")

$zj.Cells.Item(2, 1).Copy()
$zj.Cells.Item(3, 1).PasteSpecial(-4122)

$zj.Cells.Item(3, 1).Value = 1
$zj.Cells.Item(3, 2).Value = "2021-Q4"
$zj.Cells.Item(3, 3).Value = 9
$zj.Cells.Item(3, 4).Value = 3.28

$zj.Cells.Item(2, 1).Value = 0
$zj.Cells.Item(2, 2).Value = "2022-Q1"
$zj.Cells.Item(2, 3).Value = 6
$zj.Cells.Item(2, 4).Value = 0.85
